$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.856.85"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.639.31"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.34"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.520"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.94"
$ws.Range("E8").Value = "  -1.83%  "
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Value = "1.649.83"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("E14").Value = "  +3.86%  "
$ws.Range("E15").Value = "  +8.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.88"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "29.866.94"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.58"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.80"
$ws.Range("E19").Value = "  -1.12%  "
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.92"
$ws.Range("E22").Value = "  +3.30%  "
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.19"
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.73"
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.55"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").Value = "1.424.19"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").Value = "  +3.47%  "
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("E37").Value = "  -5.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0173"
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "76.92"
$ws.Range("E40").Value = "  +11.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.563"
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("E44").Value = "  -2.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.37"
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "48.98"
$ws.Range("E49").Value = "  -9.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "93.55"
$ws.Range("E50").Value = "  +6.27%  "
$ws.Range("E51").Value = "  +0.33%  "
